$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 206.66667
$ws.Range("I9").Value = 283.75
$ws.Range("K9").Value = 283.75
$ws.Range("M9").Value = -114.75
$ws.Range("H40").Value = 11210.741
$ws.Range("J40").Value = 13207.75
$ws.Range("L40").Value = 13207.75
$ws.Range("N40").Value = -13557.75
$ws.Range("H55").Value = 451.46667
$ws.Range("J55").Value = 261.5
$ws.Range("L55").Value = 261.5
$ws.Range("N55").Value = -689.5
$ws.Range("H115").Value = 594.6667
$ws.Range("I115").Value = 594.6667
$ws.Range("K115").Value = 1784.0001
$ws.Range("M115").Value = -217.0001
$ws.Range("H116").Value = 25800500
$ws.Range("I116").Value = 28515600
$ws.Range("K116").Value = 28515600
$ws.Range("M116").Value = -28512158
$ws.Range("H127").Value = 2802.0789
$ws.Range("I127").Value = 915.6667
$ws.Range("J127").Value = 3155.7812
$ws.Range("K127").Value = 2747.0001
$ws.Range("L127").Value = 9467.3436
$ws.Range("M127").Value = 2212.9999
$ws.Range("N127").Value = -19387.3436
$ws.Range("H132").Value = 8952.438
$ws.Range("I132").Value = 1317.1522
$ws.Range("K132").Value = 3951.4566
$ws.Range("M132").Value = -1421.4566
$ws.Range("H135").Value = 2446
$ws.Range("I135").Value = 889.88464
$ws.Range("K135").Value = 8008.96176
$ws.Range("M135").Value = -5473.96176
$ws.Range("H138").Value = 3330.3
$ws.Range("I138").Value = 1200.6786
$ws.Range("J138").Value = 4158.4863
$ws.Range("K138").Value = 3602.0358
$ws.Range("L138").Value = 12475.4589
$ws.Range("M138").Value = 1537.9642
$ws.Range("N138").Value = -22755.4589

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 626234.1
$ws.Range("I2").Value = 674186.75
$ws.Range("K2").Value = 674186.75
$ws.Range("M2").Value = -674073.75
$ws.Range("H28").Value = 13450
$ws.Range("I28").Value = 13450
$ws.Range("K28").Value = 13450
$ws.Range("M28").Value = -13258
$ws.Range("H32").Value = 13515.739
$ws.Range("I32").Value = 12909.431
$ws.Range("K32").Value = 12909.431
$ws.Range("M32").Value = -12622.431
$ws.Range("H41").Value = 10000
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10828
$ws.Range("H74").Value = 1584.9565
$ws.Range("I74").Value = 1175.0714
$ws.Range("K74").Value = 1175.0714
$ws.Range("M74").Value = -301.0714
$ws.Range("H77").Value = 1584.9565
$ws.Range("I77").Value = 1175.0714
$ws.Range("K77").Value = 5875.357
$ws.Range("M77").Value = -1507.357
$ws.Range("H99").Value = 13450
$ws.Range("I99").Value = 13450
$ws.Range("K99").Value = 13450
$ws.Range("M99").Value = -10455
$ws.Range("H116").Value = 626234.1
$ws.Range("I116").Value = 674186.75
$ws.Range("K116").Value = 674186.75
$ws.Range("M116").Value = -671892.75
$ws.Range("H122").Value = 5077.8237
$ws.Range("I122").Value = 2319.182
$ws.Range("J122").Value = 10135.333
$ws.Range("K122").Value = 6957.545999999999
$ws.Range("L122").Value = 30405.999
$ws.Range("M122").Value = -4507.545999999999
$ws.Range("N122").Value = -35305.999
$ws.Range("H132").Value = 13797.49
$ws.Range("I132").Value = 22762.223
$ws.Range("K132").Value = 68286.66900000001
$ws.Range("M132").Value = -65756.66900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 626234.1
$ws.Range("I3").Value = 674186.75
$ws.Range("K3").Value = 674186.75
$ws.Range("M3").Value = -674072.75
$ws.Range("H20").Value = 3162.1516
$ws.Range("I20").Value = 2883.0476
$ws.Range("K20").Value = 2883.0476
$ws.Range("M20").Value = -2636.0476
$ws.Range("H130").Value = 114999.5
$ws.Range("J130").Value = 114999.5
$ws.Range("L130").Value = 114999.5
$ws.Range("N130").Value = -125039.5
$ws.Range("H134").Value = 1138.8605
$ws.Range("I134").Value = 1108.7838
$ws.Range("K134").Value = 3326.3514
$ws.Range("M134").Value = -791.3513999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6984.643
$ws.Range("J31").Value = 7291.154
$ws.Range("L31").Value = 7291.154
$ws.Range("N31").Value = -7881.154
$ws.Range("H34").Value = 6984.643
$ws.Range("J34").Value = 7291.154
$ws.Range("L34").Value = 7291.154
$ws.Range("N34").Value = -7695.154
$ws.Range("H99").Value = 6915.375
$ws.Range("I99").Value = 4162.375
$ws.Range("K99").Value = 4162.375
$ws.Range("M99").Value = -2664.375
$ws.Range("H126").Value = 6915.375
$ws.Range("I126").Value = 4162.375
$ws.Range("K126").Value = 12487.125
$ws.Range("M126").Value = -10017.125
$ws.Range("H132").Value = 8555112
$ws.Range("I132").Value = 11915287
$ws.Range("K132").Value = 35745861
$ws.Range("M132").Value = -35743331
$ws.Range("H134").Value = 2039.4546
$ws.Range("I134").Value = 2083.5
$ws.Range("J134").Value = 1599
$ws.Range("K134").Value = 6250.5
$ws.Range("L134").Value = 4797
$ws.Range("M134").Value = -3715.5
$ws.Range("N134").Value = -9867

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6583.222
$ws.Range("J39").Value = 6788
$ws.Range("L39").Value = 20364
$ws.Range("N39").Value = -20952
$ws.Range("H55").Value = 2128.6875
$ws.Range("I55").Value = 758.1667
$ws.Range("J55").Value = 2951
$ws.Range("K55").Value = 2274.5001
$ws.Range("L55").Value = 8853
$ws.Range("M55").Value = -2097.5001
$ws.Range("N55").Value = -9207
$ws.Range("H122").Value = 1163.6923
$ws.Range("I122").Value = 765.6
$ws.Range("K122").Value = 6890.400000000001
$ws.Range("M122").Value = -4440.400000000001
$ws.Range("H127").Value = 2000
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920
$ws.Range("H131").Value = 14150716
$ws.Range("I131").Value = 953
$ws.Range("J131").Value = 16327602
$ws.Range("K131").Value = 2859
$ws.Range("L131").Value = 48982806
$ws.Range("M131").Value = 2181
$ws.Range("N131").Value = -48992886
$ws.Range("H132").Value = 5270.75
$ws.Range("I132").Value = 883
$ws.Range("K132").Value = 7947
$ws.Range("M132").Value = -5417
$ws.Range("H140").Value = 10301.333
$ws.Range("I140").Value = 4957.1665
$ws.Range("J140").Value = 17426.889
$ws.Range("K140").Value = 14871.4995
$ws.Range("L140").Value = 52280.667
$ws.Range("M140").Value = -9691.499500000002
$ws.Range("N140").Value = -62640.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1404999.9
$ws.Range("I80").Value = 2379428.8
$ws.Range("K80").Value = 2379428.8
$ws.Range("M80").Value = -2378430.8
$ws.Range("H83").Value = 1404999.9
$ws.Range("I83").Value = 2379428.8
$ws.Range("K83").Value = 11897144
$ws.Range("M83").Value = -11892152
$ws.Range("H97").Value = 693.7619
$ws.Range("I97").Value = 566.1875
$ws.Range("K97").Value = 566.1875
$ws.Range("M97").Value = -70.1875
$ws.Range("H122").Value = 737874.7
$ws.Range("I122").Value = 1377268.6
$ws.Range("J122").Value = 7138.7144
$ws.Range("K122").Value = 4131805.8
$ws.Range("L122").Value = 21416.1432
$ws.Range("M122").Value = -4129355.8
$ws.Range("N122").Value = -26316.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 950
$ws.Range("I22").Value = 957.1429000000001
$ws.Range("K22").Value = 957.1429000000001
$ws.Range("M22").Value = -662.1429000000001
$ws.Range("H27").Value = 950
$ws.Range("I27").Value = 957.1429000000001
$ws.Range("K27").Value = 957.1429000000001
$ws.Range("M27").Value = -850.1429000000001
$ws.Range("H55").Value = 306.33334
$ws.Range("J55").Value = 385.72726
$ws.Range("L55").Value = 385.72726
$ws.Range("N55").Value = -731.72726
$ws.Range("H132").Value = 5309.2856
$ws.Range("I132").Value = 5661.2856
$ws.Range("K132").Value = 16983.8568
$ws.Range("M132").Value = -14453.8568
$ws.Range("H136").Value = 2869.647
$ws.Range("I136").Value = 2002.6923
$ws.Range("K136").Value = 6008.0769
$ws.Range("M136").Value = -3458.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 841.2632
$ws.Range("I113").Value = 1088.8889
$ws.Range("J113").Value = 618.4
$ws.Range("K113").Value = 3266.6667
$ws.Range("L113").Value = 1855.2
$ws.Range("M113").Value = -1096.6667
$ws.Range("N113").Value = -6195.2
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 13514374
$ws.Range("I132").Value = 1001.5417
$ws.Range("J132").Value = 38462140
$ws.Range("K132").Value = 3004.6251
$ws.Range("L132").Value = 115386420
$ws.Range("M132").Value = -474.6251000000002
$ws.Range("N132").Value = -115391480

Write-Host "done"